$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New seed data for rows 2-11 (columns A-D). Column E ("Test Case N") is unchanged.
$data = @(
    @("Hintz LLC", 1763, "Jocelynfort", 42),
    @("Reichel, Weissnat and Stroman", 1618, "Sauerburgh", 57),
    @("Cormier - Hodkiewicz", 1763, "Mariannetown", 55),
    @("Armstrong Inc", 1809, "Rauview", 92),
    @("Padberg - Bauch", 2004, "Hillsport", 83),
    @("Gutkowski, Mraz and Moore", 1874, "Round Rock", 84),
    @("Ziemann - Lesch", 1839, "Beattyport", 21),
    @("Kuhic - Runte", 1645, "Judsonview", 46),
    @("Morar - Harvey", 1987, "Gorczanytown", 33),
    @("Glover - Rempel", 1798, "Goldenburgh", 56)
)

$row = 2
foreach ($record in $data) {
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $ws.Cells.Item($row, 4).Value = $record[3]
    $row++
}
